# Scanplan_new_EventsMP.docx edit
# -------------------------------
#  1. Move the "_GoBack" bookmark from the "Language loc / MD loc / Other
#     locs / Resting state / DTI" paragraph down to span the three
#     paragraphs "Anatomical: ... AA scout ______" through
#     "T1 (MPRAGE_1iso): ______", which is where the author's last edit
#     landed.
#  2. Bump the two "129" run counts (the BioLoc rows' last column) to
#     "149". The edit was made mid-run (changing the "2" to a "4"), so
#     each cell ends up with its text split across two runs: "14" + "9".

$d = $word.ActiveDocument

# --- 1. Relocate the "_GoBack" bookmark -----------------------------------
# Adding a bookmark named "_GoBack" automatically replaces/removes any
# existing bookmark of that name elsewhere in the document (Word only
# ever keeps a single "_GoBack" location), so the old one disappears on
# its own once we add the new one.

$findStart = $d.Content
$findStart.Find.Execute("Anatomical:", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$startPos = $findStart.Paragraphs(1).Range.Start

$findEnd = $d.Content
$findEnd.Find.Execute("MPRAGE_1iso", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$endPos = $findEnd.Paragraphs(1).Range.End

$goBackRange = $d.Range($startPos, $endPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# --- 2. "129" -> "149" wherever it shows up as a scan-count cell ---------

$search = $d.Content
while ($search.Find.Execute("129", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)) {

    $cellStart = $search.Start

    # Replace "129" with "14" in one shot ...
    $wholeRange = $d.Range($cellStart, $cellStart + 3)
    $wholeRange.Text = "14"

    # ... then append the trailing "9" as its own run, mirroring the
    # mid-edit that left the tail of the original run un-merged.
    $insertPoint = $d.Range($cellStart + 2, $cellStart + 2)
    $insertPoint.InsertAfter("9")

    $nineRange = $d.Range($cellStart + 2, $cellStart + 3)
    $nineRange.Font.Bold = 1
    $nineRange.Font.Bold = 0

    # Resume searching right after the text we just touched.
    $search = $d.Range($cellStart + 3, $d.Content.End)
}
